# Update match-odds figures per the 2025-01-28 FlashScore refresh.
# Each assignment below sets one odds cell on Sheet1 to its updated value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("G4").Value = 2.38
$ws.Range("I4").Value = 3.6
$ws.Range("J4").Value = 3.25
$ws.Range("M4").Value = 1.13
$ws.Range("N4").Value = 6
$ws.Range("U4").Value = 1.67
$ws.Range("V4").Value = 2.1
$ws.Range("W4").Value = 2.25
$ws.Range("X4").Value = 1.57
$ws.Range("Y4").Value = 5.5
$ws.Range("Z4").Value = 9.5
$ws.Range("AG4").Value = 21
$ws.Range("AJ4").Value = 7.5
$ws.Range("AL4").Value = 15
$ws.Range("AN4").Value = 41
# Row 5
$ws.Range("G5").Value = 2.88
$ws.Range("J5").Value = 3.75
$ws.Range("K5").Value = 1.83
$ws.Range("W5").Value = 2.2
$ws.Range("X5").Value = 1.62
$ws.Range("AA5").Value = 12
$ws.Range("AH5").Value = 81
# Row 7
$ws.Range("G7").Value = 1.8
$ws.Range("I7").Value = 5.25
$ws.Range("J7").Value = 2.6
$ws.Range("O7").Value = 1.57
$ws.Range("P7").Value = 2.25
$ws.Range("S7").Value = 6
$ws.Range("T7").Value = 1.13
$ws.Range("U7").Value = 1.67
$ws.Range("V7").Value = 2.1
$ws.Range("W7").Value = 2.63
$ws.Range("X7").Value = 1.44
$ws.Range("Z7").Value = 6.5
$ws.Range("AA7").Value = 10
$ws.Range("AE7").Value = 5.5
$ws.Range("AG7").Value = 26
$ws.Range("AM7").Value = 67
$ws.Range("AO7").Value = 67
$ws.Range("AR7").Value = 4.87
# Row 9
$ws.Range("M9").Value = 1.05
$ws.Range("N9").Value = 11
$ws.Range("Q9").Value = 1.88
$ws.Range("R9").Value = 1.98
$ws.Range("S9").Value = 3.25
$ws.Range("T9").Value = 1.33
$ws.Range("U9").Value = 1.4
$ws.Range("V9").Value = 2.75
$ws.Range("AD9").Value = 29
$ws.Range("AE9").Value = 11
$ws.Range("AI9").Value = 201
$ws.Range("AO9").Value = 29
# Row 10
$ws.Range("G10").Value = 2.8
$ws.Range("I10").Value = 2.3
$ws.Range("J10").Value = 3.4
$ws.Range("AF10").Value = 7
$ws.Range("AJ10").Value = 9.5
$ws.Range("AN10").Value = 17
# Row 12
$ws.Range("H12").Value = 2.9
$ws.Range("J12").Value = 2.87
$ws.Range("K12").Value = 1.91
$ws.Range("L12").Value = 4.15
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.47
$ws.Range("S12").Value = 4.45
$ws.Range("T12").Value = 1.17
$ws.Range("U12").Value = 1.53
$ws.Range("V12").Value = 2.35
$ws.Range("W12").Value = 2.05
$ws.Range("Y12").Value = 5.6
$ws.Range("Z12").Value = 9
$ws.Range("AA12").Value = 9.5
$ws.Range("AC12").Value = 22
$ws.Range("AD12").Value = 40
$ws.Range("AJ12").Value = 8.25
$ws.Range("AK12").Value = 18
$ws.Range("AL12").Value = 12.5
$ws.Range("AN12").Value = 40
$ws.Range("AO12").Value = 50
# Row 13
$ws.Range("G13").Value = 2.57
$ws.Range("H13").Value = 2.75
$ws.Range("I13").Value = 3.05
$ws.Range("J13").Value = 3.25
$ws.Range("L13").Value = 3.75
$ws.Range("P13").Value = 2.32
$ws.Range("Q13").Value = 2.57
$ws.Range("Y13").Value = 6.3
$ws.Range("Z13").Value = 11.5
$ws.Range("AB13").Value = 30
$ws.Range("AJ13").Value = 6.8
$ws.Range("AK13").Value = 14
$ws.Range("AL13").Value = 11.5
$ws.Range("AM13").Value = 40
# Row 14
$ws.Range("J14").Value = 3.25
# Row 16
$ws.Range("G16").Value = 2.45
$ws.Range("I16").Value = 2.7
$ws.Range("J16").Value = 3.1
$ws.Range("L16").Value = 3.4
$ws.Range("Y16").Value = 8.5
$ws.Range("Z16").Value = 12
$ws.Range("AA16").Value = 9.5
$ws.Range("AB16").Value = 23
$ws.Range("AC16").Value = 19
$ws.Range("AM16").Value = 29
# Row 17
$ws.Range("M17").Value = 1.05
$ws.Range("N17").Value = 11
$ws.Range("O17").Value = 1.25
$ws.Range("P17").Value = 3.75
$ws.Range("Q17").Value = 1.88
$ws.Range("R17").Value = 1.98
# Row 18
$ws.Range("N18").Value = 9
# Row 19
$ws.Range("M19").Value = 1.06
$ws.Range("N19").Value = 10
$ws.Range("Q19").Value = 2.05
$ws.Range("R19").Value = 1.75
# Row 20
$ws.Range("G20").Value = 1.5
$ws.Range("H20").Value = 4.33
$ws.Range("I20").Value = 5.5
$ws.Range("J20").Value = 2.05
$ws.Range("K20").Value = 2.4
$ws.Range("M20").Value = 1.03
$ws.Range("N20").Value = 15
$ws.Range("O20").Value = 1.2
$ws.Range("P20").Value = 4.33
$ws.Range("Q20").Value = 1.67
$ws.Range("R20").Value = 2.15
$ws.Range("Z20").Value = 7.5
$ws.Range("AA20").Value = 8.5
$ws.Range("AB20").Value = 11
$ws.Range("AE20").Value = 15
$ws.Range("AF20").Value = 8.5
$ws.Range("AI20").Value = 251
# Row 21
$ws.Range("L21").Value = 3.5
$ws.Range("Q21").Value = 1.93
$ws.Range("R21").Value = 1.93
$ws.Range("W21").Value = 1.7
$ws.Range("X21").Value = 2.05
$ws.Range("Z21").Value = 12
# Row 22
$ws.Range("G22").Value = 1.67
$ws.Range("I22").Value = 4.75
$ws.Range("M22").Value = 1.05
$ws.Range("N22").Value = 11
# Row 23
$ws.Range("H23").Value = 4.33
$ws.Range("I23").Value = 6.5
$ws.Range("J23").Value = 2
$ws.Range("M23").Value = 1.06
$ws.Range("N23").Value = 10
# Row 24
$ws.Range("M24").Value = 1.08
$ws.Range("N24").Value = 8
# Row 25
$ws.Range("G25").Value = 2.05
$ws.Range("I25").Value = 3.4
$ws.Range("W25").Value = 1.75
$ws.Range("X25").Value = 2
$ws.Range("AI25").Value = 201
$ws.Range("AK25").Value = 17
# Row 26
$ws.Range("O26").Value = 1.36
$ws.Range("P26").Value = 3
# Row 27
$ws.Range("G27").Value = 2.1
$ws.Range("H27").Value = 3.3
$ws.Range("I27").Value = 3.3
$ws.Range("L27").Value = 3.75
$ws.Range("O27").Value = 1.25
$ws.Range("P27").Value = 3.75
$ws.Range("Q27").Value = 1.85
$ws.Range("R27").Value = 1.95
$ws.Range("S27").Value = 3
$ws.Range("T27").Value = 1.36
$ws.Range("U27").Value = 1.4
$ws.Range("V27").Value = 2.75
$ws.Range("Z27").Value = 11
$ws.Range("AK27").Value = 17
$ws.Range("AM27").Value = 34
# Row 28
$ws.Range("G28").Value = 2.22
$ws.Range("H28").Value = 3.25
$ws.Range("I28").Value = 3
$ws.Range("J28").Value = 2.8
$ws.Range("K28").Value = 2.1
$ws.Range("O28").Value = 1.23
$ws.Range("P28").Value = 3.35
$ws.Range("Q28").Value = 1.7
$ws.Range("R28").Value = 1.91
$ws.Range("T28").Value = 1.38
$ws.Range("Y28").Value = 9
$ws.Range("AE28").Value = 11.25
$ws.Range("AF28").Value = 6.4
$ws.Range("AK28").Value = 17.5
$ws.Range("AO28").Value = 28
# Row 29
$ws.Range("K29").Value = 2.38
$ws.Range("N29").Value = 13
$ws.Range("Y29").Value = 7.5
$ws.Range("AE29").Value = 13
$ws.Range("AG29").Value = 17
# Row 30
$ws.Range("G30").Value = 2.5
$ws.Range("I30").Value = 2.9
$ws.Range("Q30").Value = 2.05
$ws.Range("R30").Value = 1.8
$ws.Range("S30").Value = 3.5
$ws.Range("T30").Value = 1.29
$ws.Range("U30").Value = 1.44
$ws.Range("V30").Value = 2.63
$ws.Range("AA30").Value = 10
$ws.Range("AH30").Value = 51
$ws.Range("AI30").Value = 251
$ws.Range("AJ30").Value = 9
# Row 31
$ws.Range("AP31").Value = 1.83
$ws.Range("AQ31").Value = 1.98
# Row 32
$ws.Range("U32").Value = 1.57
$ws.Range("V32").Value = 2.25
$ws.Range("Y32").Value = 6.5
$ws.Range("AP32").Value = 1.85
$ws.Range("AQ32").Value = 1.95
# Row 33
$ws.Range("G33").Value = 2.63
$ws.Range("Q33").Value = 2
$ws.Range("R33").Value = 1.85
